# Code réécrit et avec commentaire
# Il reste la gestion de s'il n'y a pas de port de connecter à s'occuper

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The day-header row only goes from 1 to 31 (one column per day of month),
# but this sheet was generated for a short month - drop the unused trailing
# day columns (N:AG, i.e. days 12-31) entirely so the sheet dimension
# shrinks back down to A1:M2.
$ws.Range("N1:AG1").EntireColumn.Delete()

# Day 9 (column K) now has a recorded value of 0.88 instead of the 0.00
# placeholder. Force the cell to stay text (it already was "0.00\n") so we
# don't turn it into a numeric cell, then restore the default "Normal"
# style so we don't leave a stray number-format behind.
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "0.88`n"
$ws.Range("K2").Style = "Normal"

# Day 11 (column M, the new last column after the delete above) has no
# recorded value yet - clear it out to an empty text cell (instead of
# deleting it outright) and reset its style back to default.
$ws.Range("M2").Value = "'"
$ws.Range("M2").Style = "Normal"
